# "Update hosts, add IgY"
#
# 1. Terminology sheet: rewrite the Host column (A) with the new, more
#    descriptive species names and add two new hosts (mouse, alpaca).
#    Also add the missing "IgY" isotype in column B.
# 2. Antibodies sheet: update the existing Host values to match the new
#    descriptive naming, and widen the two list-validation ranges that
#    point at the Terminology sheet to cover the newly added rows.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Terminology sheet (protected) - unprotect, edit, re-protect so the
# sheet ends up locked again exactly like before.
# ------------------------------------------------------------------
$term = $wb.Worksheets.Item("Terminology")
$term.Unprotect()

$term.Range("A2").Value = "chicken (Gallus gallus)"
$term.Range("A3").Value = "human (Homo sapiens)"
$term.Range("A4").Value = "llama (Lama glama)"
$term.Range("A5").Value = "mouse (Mus musculus)"
$term.Range("A6").Value = "alpaca (Vicugna pacos)"

$term.Range("B16").Value = "IgY"

$term.Protect()

# ------------------------------------------------------------------
# Antibodies sheet - update Host column values to the new names, and
# extend the data-validation source ranges on the Terminology sheet to
# include the newly added rows (hosts A2:A6, isotypes B2:B16).
# ------------------------------------------------------------------
$ab = $wb.Worksheets.Item("Antibodies")

$ab.Range("B2").Value = "human (Homo sapiens)"
$ab.Range("B3").Value = "human (Homo sapiens)"
$ab.Range("B4").Value = "mouse (Mus musculus)"
$ab.Range("B5").Value = "human (Homo sapiens)"
$ab.Range("B6").Value = "mouse (Mus musculus)"
$ab.Range("B7").Value = "mouse (Mus musculus)"
$ab.Range("B8").Value = "human (Homo sapiens)"
$ab.Range("B9").Value = "mouse (Mus musculus)"
$ab.Range("B10").Value = "human (Homo sapiens)"
$ab.Range("B11").Value = "mouse (Mus musculus)"

$ab.Range("B2:B100").Validation.Formula1 = "=Terminology!`$A`$2:`$A`$6"
$ab.Range("C2:C100").Validation.Formula1 = "=Terminology!`$B`$2:`$B`$16"
